$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 1980.1837
$ws.Range("J138").Value = 2159.8948
$ws.Range("L138").Value = 6479.6844
$ws.Range("N138").Value = -16759.6844

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2190.04
$ws.Range("I2").Value = 1696.1666
$ws.Range("J2").Value = 2645.923
$ws.Range("K2").Value = 1696.1666
$ws.Range("L2").Value = 2645.923
$ws.Range("M2").Value = -1583.1666
$ws.Range("N2").Value = -2871.923
$ws.Range("H50").Value = 1868.5
$ws.Range("J50").Value = 2533.3333
$ws.Range("L50").Value = 2533.3333
$ws.Range("N50").Value = -3961.3333
$ws.Range("H74").Value = 2350.6875
$ws.Range("I74").Value = 2058.0833
$ws.Range("J74").Value = 3228.5
$ws.Range("K74").Value = 2058.0833
$ws.Range("L74").Value = 3228.5
$ws.Range("M74").Value = -1184.0833
$ws.Range("N74").Value = -4976.5
$ws.Range("H77").Value = 2350.6875
$ws.Range("I77").Value = 2058.0833
$ws.Range("J77").Value = 3228.5
$ws.Range("K77").Value = 10290.4165
$ws.Range("L77").Value = 16142.5
$ws.Range("M77").Value = -5922.416499999999
$ws.Range("N77").Value = -24878.5
$ws.Range("H97").Value = 391.5
$ws.Range("J97").Value = 561.6667
$ws.Range("L97").Value = 561.6667
$ws.Range("N97").Value = -1553.6667
$ws.Range("H110").Value = 3336.5715
$ws.Range("I110").Value = 3352.182
$ws.Range("J110").Value = 3279.3333
$ws.Range("K110").Value = 3352.182
$ws.Range("L110").Value = 3279.3333
$ws.Range("M110").Value = -1307.182
$ws.Range("N110").Value = -7369.3333
$ws.Range("H116").Value = 2190.04
$ws.Range("I116").Value = 1696.1666
$ws.Range("J116").Value = 2645.923
$ws.Range("K116").Value = 1696.1666
$ws.Range("L116").Value = 2645.923
$ws.Range("M116").Value = 597.8334
$ws.Range("N116").Value = -7233.923
$ws.Range("H122").Value = 4003.724
$ws.Range("I122").Value = 3884.6
$ws.Range("J122").Value = 4748.25
$ws.Range("K122").Value = 11653.8
$ws.Range("L122").Value = 14244.75
$ws.Range("M122").Value = -9203.799999999999
$ws.Range("N122").Value = -19144.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2190.04
$ws.Range("I3").Value = 1696.1666
$ws.Range("J3").Value = 2645.923
$ws.Range("K3").Value = 1696.1666
$ws.Range("L3").Value = 2645.923
$ws.Range("M3").Value = -1582.1666
$ws.Range("N3").Value = -2873.923
$ws.Range("H41").Value = 0
$ws.Range("J41").Value = 0
$ws.Range("N41").Value = 0
$ws.Range("H42").Value = 110000
$ws.Range("J42").Value = 110000
$ws.Range("L42").Value = 110000
$ws.Range("N42").Value = -110656
$ws.Range("H46").Value = 15210.526
$ws.Range("H48").Value = 0
$ws.Range("J48").Value = 0
$ws.Range("N48").Value = 0
$ws.Range("H64").Value = 1125.1818
$ws.Range("J64").Value = 897.125
$ws.Range("L64").Value = 897.125
$ws.Range("N64").Value = -1347.125
$ws.Range("H67").Value = 1125.1818
$ws.Range("J67").Value = 897.125
$ws.Range("L67").Value = 897.125
$ws.Range("N67").Value = -2457.125
$ws.Range("L41").ClearContents()
$ws.Range("L48").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2567.8823
$ws.Range("I31").Value = 3284.1
$ws.Range("J31").Value = 1544.7142
$ws.Range("K31").Value = 3284.1
$ws.Range("L31").Value = 1544.7142
$ws.Range("M31").Value = -2989.1
$ws.Range("N31").Value = -2134.7142
$ws.Range("H34").Value = 2567.8823
$ws.Range("I34").Value = 3284.1
$ws.Range("J34").Value = 1544.7142
$ws.Range("K34").Value = 3284.1
$ws.Range("L34").Value = 1544.7142
$ws.Range("M34").Value = -3082.1
$ws.Range("N34").Value = -1948.7142
$ws.Range("H98").Value = 26332.666
$ws.Range("J98").Value = 26332.666
$ws.Range("L98").Value = 26332.666
$ws.Range("N98").Value = -30824.666
$ws.Range("H99").Value = 59261936
$ws.Range("I99").Value = 33336130
$ws.Range("J99").Value = 111113550
$ws.Range("K99").Value = 33336130
$ws.Range("L99").Value = 111113550
$ws.Range("M99").Value = -33334632
$ws.Range("N99").Value = -111116546
$ws.Range("H126").Value = 59261936
$ws.Range("I126").Value = 33336130
$ws.Range("J126").Value = 111113550
$ws.Range("K126").Value = 100008390
$ws.Range("L126").Value = 333340650
$ws.Range("M126").Value = -100005920
$ws.Range("N126").Value = -333345590

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H33").Value = 466.75
$ws.Range("J33").Value = 506.66666
$ws.Range("L33").Value = 3039.99996
$ws.Range("N33").Value = -3605.99996
$ws.Range("H41").Value = 3376.1428
$ws.Range("I41").Value = 2999
$ws.Range("J41").Value = 3439
$ws.Range("K41").Value = 8997
$ws.Range("L41").Value = 10317
$ws.Range("M41").Value = -8659
$ws.Range("N41").Value = -10993
$ws.Range("H59").Value = 100
$ws.Range("I59").Value = 0
$ws.Range("K59").Value = 0
$ws.Range("H62").Value = 1729.75
$ws.Range("I62").Value = 973.3333
$ws.Range("J62").Value = 3999
$ws.Range("K62").Value = 2919.9999
$ws.Range("L62").Value = 11997
$ws.Range("M62").Value = -2233.9999
$ws.Range("N62").Value = -13369
$ws.Range("H65").Value = 1729.75
$ws.Range("I65").Value = 973.3333
$ws.Range("J65").Value = 3999
$ws.Range("K65").Value = 8759.9997
$ws.Range("L65").Value = 35991
$ws.Range("M65").Value = -5327.9997
$ws.Range("N65").Value = -42855
$ws.Range("H69").Value = 6300
$ws.Range("J69").Value = 1500
$ws.Range("L69").Value = 4500
$ws.Range("N69").Value = -6122
$ws.Range("H72").Value = 6300
$ws.Range("J72").Value = 1500
$ws.Range("L72").Value = 13500
$ws.Range("N72").Value = -21612
$ws.Range("H131").Value = 3166
$ws.Range("I131").Value = 2249.25
$ws.Range("J131").Value = 4999.5
$ws.Range("K131").Value = 6747.75
$ws.Range("L131").Value = 14998.5
$ws.Range("M131").Value = -1707.75
$ws.Range("N131").Value = -25078.5
$ws.Range("H132").Value = 0
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 0
$ws.Range("N132").Value = 0
$ws.Range("M59").ClearContents()
$ws.Range("L132").ClearContents()
$ws.Range("M132").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 4559.65
$ws.Range("J102").Value = 5000
$ws.Range("L102").Value = 5000
$ws.Range("N102").Value = -8244
$ws.Range("H113").Value = 1841.9
$ws.Range("J113").Value = 1987.1428
$ws.Range("L113").Value = 1987.1428
$ws.Range("N113").Value = -6327.1428
$ws.Range("H126").Value = 5959
$ws.Range("J126").Value = 10000
$ws.Range("L126").Value = 30000
$ws.Range("N126").Value = -34940
$ws.Range("H132").Value = 9256.571
$ws.Range("I132").Value = 9132.666999999999
$ws.Range("K132").Value = 27398.001
$ws.Range("M132").Value = -24868.001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 1179.8
$ws.Range("I61").Value = 1179.8
$ws.Range("K61").Value = 1179.8
$ws.Range("M61").Value = -977.8
$ws.Range("H93").Value = 966.7646999999999
$ws.Range("I93").Value = 805
$ws.Range("K93").Value = 805
$ws.Range("M93").Value = 443
$ws.Range("H113").Value = 1179.8
$ws.Range("I113").Value = 1179.8
$ws.Range("K113").Value = 1179.8
$ws.Range("M113").Value = 990.2
$ws.Range("H132").Value = 3502.842
$ws.Range("J132").Value = 4999.5
$ws.Range("L132").Value = 14998.5
$ws.Range("N132").Value = -20058.5
